$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated reaction-coverage data (rows 5-33, columns B-F)
$data = @(
    @(5, 62, 65.113182423435418, 66.577896138482018, 65.113182423435418, 66.577896138482018),
    @(6, 62.8, 66.178428761651134, 67.243675099866849, 66.178428761651134, 67.243675099866849),
    @(7, 45.866666666666667, 45.80559254327563, 45.80559254327563, 45.80559254327563, 45.80559254327563),
    @(8, 52.533333333333331, 52.996005326231689, 54.327563249001329, 52.996005326231689, 54.327563249001329),
    @(9, 52, 53.129161118508655, 54.593874833555262, 53.129161118508655, 54.593874833555262),
    @(10, 49.866666666666667, 49.800266311584551, 49.800266311584551, 49.800266311584551, 49.800266311584551),
    @(11, 62.133333333333333, 62.583222370173104, 66.045272969374167, 62.583222370173104, 66.045272969374167),
    @(12, 62.133333333333333, 66.577896138482018, 71.770972037283627, 66.577896138482018, 71.770972037283627),
    @(13, 30.8, 33.022636484687084, 34.886817576564582, 35.552596537949398, 40.213049267643143),
    @(14, 32.133333333333333, 34.221038615179758, 36.085219707057256, 36.75099866844208, 41.411451398135817),
    @(15, 29.066666666666666, 30.359520639147803, 32.223701731025301, 31.824234354194406, 36.484687083888147),
    @(16, 28.4, 29.693741677762983, 31.424766977363515, 31.158455392809586, 35.552596537949398),
    @(17, 30.8, 33.022636484687084, 34.886817576564582, 35.552596537949398, 40.213049267643143),
    @(18, 38.533333333333331, 41.944074567243675, 50.599201065246341, 44.873501997336881, 55.126498002663119),
    @(19, 37.6, 40.612516644474034, 49.134487350199734, 43.54194407456724, 53.661784287616513),
    @(20, 33.466666666666669, 22.63648468708389, 37.949400798934754, 24.900133155792275, 43.408788282290281),
    @(21, 34, 22.902796271637815, 38.482023968042611, 25.166444740346204, 43.941411451398139),
    @(22, 33.466666666666669, 22.769640479360852, 38.08255659121172, 25.033288948069242, 43.54194407456724),
    @(23, 36.4, 24.500665778961384, 41.278295605858858, 26.897470039946736, 47.270306258322236),
    @(24, 36.666666666666664, 24.63382157123835, 41.67776298268975, 27.030625832223702, 47.669773635153128),
    @(25, 36.666666666666664, 24.766977363515313, 41.810918774966709, 27.163781624500665, 47.802929427430094),
    @(26, 24.533333333333335, 24.63382157123835, 25.29960053262317, 26.631158455392811, 29.826897470039945),
    @(27, 25.066666666666666, 25.166444740346204, 25.832223701731024, 27.296937416777631, 30.492676431424766),
    @(28, 24.133333333333333, 24.367509986684421, 24.766977363515313, 26.364846870838882, 29.294274300932091),
    @(29, 21.6, 21.837549933422103, 22.370173102529961, 23.701731025299601, 26.364846870838882),
    @(30, 24.8, 24.900133155792275, 25.565912117177096, 26.897470039946736, 30.093209054593874),
    @(31, 36.133333333333333, 24.101198402130493, 41.278295605858858, 26.498002663115845, 47.403462050599202),
    @(32, 34.93333333333333, 23.834886817576564, 39.280958721704394, 25.699067909454062, 44.474034620505989),
    @(33, 34.93333333333333, 24.234354194407455, 39.280958721704394, 26.231691078561916, 44.074567243675098)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("B$r").Value = $entry[1]
    $ws.Range("C$r").Value = $entry[2]
    $ws.Range("D$r").Value = $entry[3]
    $ws.Range("E$r").Value = $entry[4]
    $ws.Range("F$r").Value = $entry[5]
}

# Column B now carries its own distinct number format (0.0) separate from
# columns C:F (0.00); reset to the base style first so the new formats are
# written cleanly (no inherited center alignment).
$ws.Range("B5:F33").Style = "Normal"
$ws.Range("C5:F33").NumberFormat = "0.00"
$ws.Range("B5:B33").NumberFormat = "0.0"

# Restore the selection reflected in the saved workbook view
$ws.Range("C5:F33").Select()
